$wb = $excel.ActiveWorkbook

# The 77decf4e-3ffd-46ea-8bff-6d401b1fb8cf.md file has finished its
# handback to both locales - update the status + handback timestamps
# across the Overview and per-locale report sheets.

$status = "Handed back: in sync with en-US"

# --- Overview sheet: row for 77decf4e...md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# --- zh-cn sheet: row for 77decf4e...md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-03-12 12:36:27"

# --- de-de sheet: row for 77decf4e...md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-03-12 12:36:32"
